$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sD2 = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.825.74"
$ws.Range("D2").Style = $sD2
$ws.Range("E2").Value = "  -2.82%  "

$sD3 = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.790.82"
$ws.Range("D3").Style = $sD3
$ws.Range("E3").Value = "  -0.57%  "

$sD4 = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = $sD4
$ws.Range("E4").Value = "  -0.14%  "

$sD5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.24"
$ws.Range("D5").Style = $sD5
$ws.Range("E5").Value = "  +0.06%  "

$ws.Range("E6").Value = "  -0.18%  "

$sD7 = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5329"
$ws.Range("D7").Style = $sD7
$ws.Range("E7").Value = "  +0.39%  "

$sD8 = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3835"
$ws.Range("D8").Style = $sD8
$ws.Range("E8").Value = "  +1.78%  "

$sD9 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07416"
$ws.Range("D9").Style = $sD9
$ws.Range("E9").Value = "  -1.02%  "

$sD10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.29"
$ws.Range("D10").Style = $sD10
$ws.Range("E10").Value = "  -2.83%  "

$sD11 = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.083"
$ws.Range("D11").Style = $sD11
$ws.Range("E11").Value = "  -2.66%  "

$sD12 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = $sD12
$ws.Range("E12").Value = "  -0.18%  "

$sD13 = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.174"
$ws.Range("D13").Style = $sD13
$ws.Range("E13").Value = "  +0.45%  "

$sD14 = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.432"
$ws.Range("D14").Style = $sD14
$ws.Range("E14").Value = "  +1.26%  "

$sD15 = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.26"
$ws.Range("D15").Style = $sD15
$ws.Range("E15").Value = "  -2.05%  "

$sD16 = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.792.88"
$ws.Range("D16").Style = $sD16
$ws.Range("E16").Value = "  -0.38%  "

$sD17 = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.99"
$ws.Range("D17").Style = $sD17
$ws.Range("E17").Value = "  -2.52%  "

$sD18 = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001056"
$ws.Range("D18").Style = $sD18
$ws.Range("E18").Value = "  -0.68%  "

$sD19 = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06519"
$ws.Range("D19").Style = $sD19
$ws.Range("E19").Value = "  +0.80%  "

$sD20 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9993"
$ws.Range("D20").Style = $sD20
$ws.Range("E20").Value = "  -0.24%  "

$sD21 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.19"
$ws.Range("D21").Style = $sD21
$ws.Range("E21").Value = "  -0.11%  "

$sD22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.942"
$ws.Range("D22").Style = $sD22
$ws.Range("E22").Value = "  +0.78%  "

$sD23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.869.28"
$ws.Range("D23").Style = $sD23
$ws.Range("E23").Value = "  -2.62%  "

$sD24 = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.11"
$ws.Range("D24").Style = $sD24
$ws.Range("E24").Value = "  +0.30%  "

$sD25 = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.094"
$ws.Range("D25").Style = $sD25
$ws.Range("E25").Value = "  -0.02%  "

$sD26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.91"
$ws.Range("D26").Style = $sD26
$ws.Range("E26").Value = "  -1.72%  "

$sD27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.07"
$ws.Range("D27").Style = $sD27
$ws.Range("E27").Value = "  -1.69%  "

$sD28 = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.997.95"
$ws.Range("D28").Style = $sD28
$ws.Range("E28").Value = "  -0.10%  "

$sD29 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.299"
$ws.Range("D29").Style = $sD29
$ws.Range("E29").Value = "  -1.81%  "

$sD30 = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.11"
$ws.Range("D30").Style = $sD30
$ws.Range("E30").Value = "  -1.09%  "

$sD31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1091"
$ws.Range("D31").Style = $sD31
$ws.Range("E31").Value = "  +3.75%  "

$sD32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.096"
$ws.Range("D32").Style = $sD32
$ws.Range("E32").Value = "  -0.49%  "

$sD33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.651"
$ws.Range("D33").Style = $sD33
$ws.Range("E33").Value = "  -1.10%  "

$sD34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.486"
$ws.Range("D34").Style = $sD34
$ws.Range("E34").Value = "  -2.66%  "

$sD35 = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06924"
$ws.Range("D35").Style = $sD35
$ws.Range("E35").Value = "  +8.10%  "

$sD36 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2194"
$ws.Range("D36").Style = $sD36
$ws.Range("E36").Value = "  -2.21%  "

$sD37 = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02262"
$ws.Range("D37").Style = $sD37
$ws.Range("E37").Value = "  -1.92%  "

$sD38 = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.030"
$ws.Range("D38").Style = $sD38
$ws.Range("E38").Value = "  -0.16%  "

$sD39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.26"
$ws.Range("D39").Style = $sD39
$ws.Range("E39").Value = "  +0.21%  "

$sD40 = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.350"
$ws.Range("D40").Style = $sD40
$ws.Range("E40").Value = "  -5.14%  "

$sD41 = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6084"
$ws.Range("D41").Style = $sD41
$ws.Range("E41").Value = "  -1.84%  "

$sD42 = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.173"
$ws.Range("D42").Style = $sD42
$ws.Range("E42").Value = "  -3.18%  "

$sD43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.411"
$ws.Range("D43").Style = $sD43
$ws.Range("E43").Value = "  -0.02%  "

$sD44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.22"
$ws.Range("D44").Style = $sD44
$ws.Range("E44").Value = "  -0.53%  "

$sD45 = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.681"
$ws.Range("D45").Style = $sD45
$ws.Range("E45").Value = "  -0.22%  "

$sD46 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5682"
$ws.Range("D46").Style = $sD46
$ws.Range("E46").Value = "  -2.75%  "

$sD47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "124.06"
$ws.Range("D47").Style = $sD47
$ws.Range("E47").Value = "  -1.27%  "

$sD48 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.908"
$ws.Range("D48").Style = $sD48
$ws.Range("E48").Value = "  -1.55%  "

$sD49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.167"
$ws.Range("D49").Style = $sD49
$ws.Range("E49").Value = "  +1.63%  "

$sD50 = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06791"
$ws.Range("D50").Style = $sD50
$ws.Range("E50").Value = "  -1.42%  "

$sD51 = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00000000295"
$ws.Range("D51").Style = $sD51
$ws.Range("E51").Value = "  +38.60%  "
